# Update average_county_temperature (col K) with refreshed NOAA temperature
# data, and recompute the dependent worst/best ASHP COP columns (R, S) for
# the affected rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 12.93898809523811
$ws.Range("R2").Value = 1.586442583591966
$ws.Range("S2").Value = 1.700608911205746

# Row 3
$ws.Range("K3").Value = 12.93898809523811

# Row 7
$ws.Range("K7").Value = -0.763888888888889
$ws.Range("R7").Value = 1.432007893438579
$ws.Range("S7").Value = 1.522400424853956

# Row 11
$ws.Range("K11").Value = 19.65277777777778
$ws.Range("R11").Value = 1.674945181765724
$ws.Range("S11").Value = 1.804078036500944

# Row 14
$ws.Range("K14").Value = 5.462962962962945
$ws.Range("R14").Value = 1.49828630419821
$ws.Range("S14").Value = 1.598520446096654

# Row 16
$ws.Range("K16").Value = 5.462962962962945
$ws.Range("R16").Value = 1.49828630419821
$ws.Range("S16").Value = 1.598520446096654

# Row 18
$ws.Range("K18").Value = 5.462962962962945
$ws.Range("R18").Value = 1.49828630419821
$ws.Range("S18").Value = 1.598520446096654

# Row 21
$ws.Range("K21").Value = 1.791666666666668
$ws.Range("R21").Value = 1.458486584262888
$ws.Range("S21").Value = 1.552746181345467

# Row 22
$ws.Range("K22").Value = 1.791666666666668

# Row 24
$ws.Range("K24").Value = 12.66820987654322
$ws.Range("R24").Value = 1.583068924143447
$ws.Range("S24").Value = 1.696684247214952

# Row 26
$ws.Range("K26").Value = 5.462962962962945
$ws.Range("R26").Value = 1.49828630419821
$ws.Range("S26").Value = 1.598520446096654

# Row 27
$ws.Range("K27").Value = 19.60879629629628
$ws.Range("R27").Value = 1.674333288469303
$ws.Range("S27").Value = 1.803359265239363

# Row 28
$ws.Range("K28").Value = 14.96875
$ws.Range("R28").Value = 1.612196950762309
$ws.Range("S28").Value = 1.730616680249932

# Row 29
$ws.Range("K29").Value = 14.96875

# Row 30
$ws.Range("K30").Value = 15.36574074074072
$ws.Range("R30").Value = 1.617332194197838
$ws.Range("S30").Value = 1.73660999151892

# Row 31
$ws.Range("K31").Value = 13.46442495126706
$ws.Range("R31").Value = 1.593030259848797
$ws.Range("S31").Value = 1.708276634982499

# Row 37
$ws.Range("K37").Value = 19.65277777777778
$ws.Range("R37").Value = 1.674945181765724
$ws.Range("S37").Value = 1.804078036500944

# Row 38
$ws.Range("K38").Value = 14.47727272727272
$ws.Range("R38").Value = 1.605884483070795
$ws.Range("S38").Value = 1.723253983867794

# Row 39
$ws.Range("K39").Value = 14.47727272727272

# Row 40
$ws.Range("K40").Value = 19.65277777777778
$ws.Range("R40").Value = 1.674945181765724
$ws.Range("S40").Value = 1.804078036500944

# Row 42
$ws.Range("K42").Value = 12.93898809523811
$ws.Range("R42").Value = 1.586442583591966
$ws.Range("S42").Value = 1.700608911205746

# Row 43
$ws.Range("K43").Value = 1.925925925925943
$ws.Range("R43").Value = 1.459904774678112
$ws.Range("S43").Value = 1.554373915558126

Write-Host "Updated average_county_temperature and ASHP COP columns for affected rows."
